{"js": "// This document is a daily \"three-digit \u00d7 one-digit multiplication\"\n// worksheet. The heading paragraph holds a date string, and the table\n// cells hold \"A\u00d7B=C\" equation strings. The edit updates the date and\n// swaps in a new set of equations, each old value being unique in the\n// document, so a direct search-and-replace per pair is safe and exact.\nconst replacements = [\n  [\"2026-02-14 Saturday\", \"2026-02-15 Sunday\"],\n  [\"296\u00d79=2664\", \"840\u00d76=5040\"],\n  [\"355\u00d79=3195\", \"434\u00d72=868\"],\n  [\"449\u00d77=3143\", \"363\u00d77=2541\"],\n  [\"766\u00d77=5362\", \"298\u00d78=2384\"],\n  [\"999\u00d75=4995\", \"880\u00d74=3520\"],\n  [\"812\u00d77=5684\", \"397\u00d76=2382\"],\n  [\"951\u00d76=5706\", \"730\u00d79=6570\"],\n  [\"826\u00d76=4956\", \"515\u00d78=4120\"],\n  [\"922\u00d75=4610\", \"378\u00d77=2646\"],\n  [\"926\u00d77=6482\", \"420\u00d75=2100\"],\n  [\"847\u00d74=3388\", \"457\u00d73=1371\"],\n  [\"423\u00d77=2961\", \"710\u00d77=4970\"],\n  [\"990\u00d72=1980\", \"873\u00d77=6111\"],\n  [\"485\u00d72=970\", \"357\u00d72=714\"],\n  [\"695\u00d79=6255\", \"380\u00d78=3040\"],\n  [\"703\u00d77=4921\", \"799\u00d76=4794\"],\n  [\"230\u00d72=460\", \"721\u00d77=5047\"],\n  [\"336\u00d74=1344\", \"242\u00d72=484\"],\n  [\"207\u00d74=828\", \"423\u00d73=1269\"],\n  [\"884\u00d78=7072\", \"494\u00d76=2964\"],\n  [\"520\u00d74=2080\", \"526\u00d77=3682\"],\n  [\"478\u00d74=1912\", \"726\u00d77=5082\"],\n  [\"369\u00d79=3321\", \"254\u00d75=1270\"],\n  [\"935\u00d79=8415\", \"431\u00d76=2586\"],\n  [\"409\u00d79=3681\", \"978\u00d76=5868\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# This document is a daily \"three-digit x one-digit multiplication\"\n# worksheet. The heading paragraph holds a date string, and the table\n# cells hold \"A x B=C\" equation strings. The edit updates the date and\n# swaps in a new set of equations. Each old value is unique in the\n# document, so a direct Find/Replace per pair is safe and exact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-02-14 Saturday\", \"2026-02-15 Sunday\"),\n    @(\"296\u00d79=2664\", \"840\u00d76=5040\"),\n    @(\"355\u00d79=3195\", \"434\u00d72=868\"),\n    @(\"449\u00d77=3143\", \"363\u00d77=2541\"),\n    @(\"766\u00d77=5362\", \"298\u00d78=2384\"),\n    @(\"999\u00d75=4995\", \"880\u00d74=3520\"),\n    @(\"812\u00d77=5684\", \"397\u00d76=2382\"),\n    @(\"951\u00d76=5706\", \"730\u00d79=6570\"),\n    @(\"826\u00d76=4956\", \"515\u00d78=4120\"),\n    @(\"922\u00d75=4610\", \"378\u00d77=2646\"),\n    @(\"926\u00d77=6482\", \"420\u00d75=2100\"),\n    @(\"847\u00d74=3388\", \"457\u00d73=1371\"),\n    @(\"423\u00d77=2961\", \"710\u00d77=4970\"),\n    @(\"990\u00d72=1980\", \"873\u00d77=6111\"),\n    @(\"485\u00d72=970\", \"357\u00d72=714\"),\n    @(\"695\u00d79=6255\", \"380\u00d78=3040\"),\n    @(\"703\u00d77=4921\", \"799\u00d76=4794\"),\n    @(\"230\u00d72=460\", \"721\u00d77=5047\"),\n    @(\"336\u00d74=1344\", \"242\u00d72=484\"),\n    @(\"207\u00d74=828\", \"423\u00d73=1269\"),\n    @(\"884\u00d78=7072\", \"494\u00d76=2964\"),\n    @(\"520\u00d74=2080\", \"526\u00d77=3682\"),\n    @(\"478\u00d74=1912\", \"726\u00d77=5082\"),\n    @(\"369\u00d79=3321\", \"254\u00d75=1270\"),\n    @(\"935\u00d79=8415\", \"431\u00d76=2586\"),\n    @(\"409\u00d79=3681\", \"978\u00d76=5868\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
